$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.353.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5352"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2667"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06409"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07854"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.667.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.895.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5543"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8204"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.372.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.696"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.055"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.011"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1233"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.220"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05879"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.644"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.291"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9726"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.835"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01603"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8686"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.065.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.847"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.807.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4388"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.014"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05168"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "
